$d = $word.ActiveDocument

# Locate the last paragraph in the document body (the empty paragraph
# immediately before the final section properties) and collapse to its end.
$lastIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastIndex)
$r = $p.Range
$r.Collapse(0)

# New paragraph 1: empty "Normal" paragraph (blank spacer line).
$r.InsertParagraphAfter()

# New paragraph 2: another empty "Normal" paragraph (blank spacer line).
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# New paragraph 3: "General:" heading line.
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter("General:")

# New paragraph 4: empty paragraph with 9pt (sz=18 half-points) font.
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r5 = $p5.Range
$r5.Font.Size = 9
$r5.Font.SizeBi = 9

# New paragraph 5: description text, also at 9pt.
$r5b = $p5.Range
$r5b.Collapse(0)
$r5b.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r6 = $p6.Range
$r6.Collapse(0)
$r6.InsertAfter("Over the whole employment process the information traffic between company and applicants (invitations, rejections, etc.) will get automated/standardized. Forms will be defined to enable standardized communication.")
$r6full = $p6.Range
$r6full.Font.Size = 9
$r6full.Font.SizeBi = 9
